$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Parameters sheet: environment value changed from "remote" to "local"
# ------------------------------------------------------------------
$wsParams = $wb.Worksheets.Item("Parameters")
$wsParams.Range("B2").Value = "local"

# ------------------------------------------------------------------
# Scenarios sheet: append five new rows for the new
# ops.web.tests.VisitDetailsModalTest scenarios.
# ------------------------------------------------------------------
$wsScenarios = $wb.Worksheets.Item("Scenarios")

$newRows = @(
    @{ Row = 4; TestCase = "changeProviderManualTimeSet" },
    @{ Row = 5; TestCase = "startVisit" },
    @{ Row = 6; TestCase = "endVisit" },
    @{ Row = 7; TestCase = "editVisitSymptoms" },
    @{ Row = 8; TestCase = "refundVisitTotal" }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $wsScenarios.Range("A$r").Value = "Y"
    $wsScenarios.Range("B$r").Value = "ops.web.tests.VisitDetailsModalTest"
    $wsScenarios.Range("D$r").Value = $entry.TestCase
}

# ------------------------------------------------------------------
# Selection / active sheet: Scenarios becomes the active tab (moving
# tabSelected away from Parameters), with D8 - the last new row's
# Test Case cell - selected.
# ------------------------------------------------------------------
[void]$wsScenarios.Activate()
[void]$wsScenarios.Range("D8").Select()

# ------------------------------------------------------------------
# Workbook-level absolute path correction (local dev folder rename,
# reflected in the x15ac:absPath metadata written at save time).
# ------------------------------------------------------------------
$wb.Path = "/Users/vahanmelikyan/Documents/heal/qa-automation/Automation/runs/"
